# Auto update Excel log
$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# --- Proximity sheet: append new EXIT row (row 22) ---
$wsProximity = $wb.Worksheets.Item("Proximity")
Set-TextCell $wsProximity "A22" "2026-02-01"
Set-TextCell $wsProximity "B22" "18:16:07"
Set-TextCell $wsProximity "C22" "18:00"
Set-TextCell $wsProximity "D22" "Living Room Main Door"
Set-TextCell $wsProximity "E22" "EXIT"
Set-TextCell $wsProximity "F22" "User EXITED Living Room Main Door"

# --- Camera sheet: append two new "Image Captured" rows (rows 21-22) ---
$wsCamera = $wb.Worksheets.Item("Camera")

Set-TextCell $wsCamera "A21" "2026-02-01"
Set-TextCell $wsCamera "B21" "18:15:55"
Set-TextCell $wsCamera "C21" "18:00"
Set-TextCell $wsCamera "D21" "Living Room Main Door"
Set-TextCell $wsCamera "E21" "Image Captured"
Set-TextCell $wsCamera "F21" "Active"

Set-TextCell $wsCamera "A22" "2026-02-01"
Set-TextCell $wsCamera "B22" "18:16:09"
Set-TextCell $wsCamera "C22" "18:00"
Set-TextCell $wsCamera "D22" "Living Room Main Door"
Set-TextCell $wsCamera "E22" "Image Captured"
Set-TextCell $wsCamera "F22" "Active"
